$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 640.05
$ws.Range("I33").Value = 533.5333000000001
$ws.Range("J33").Value = 959.6
$ws.Range("K33").Value = 533.5333000000001
$ws.Range("L33").Value = 959.6
$ws.Range("M33").Value = -304.5333000000001
$ws.Range("N33").Value = -1417.6
$ws.Range("H64").Value = 55562076
$ws.Range("I64").Value = 166668480
$ws.Range("J64").Value = 8875.5
$ws.Range("K64").Value = 166668480
$ws.Range("L64").Value = 8875.5
$ws.Range("M64").Value = -166668232
$ws.Range("N64").Value = -9371.5
$ws.Range("H67").Value = 55562076
$ws.Range("I67").Value = 166668480
$ws.Range("J67").Value = 8875.5
$ws.Range("K67").Value = 166668480
$ws.Range("L67").Value = 8875.5
$ws.Range("M67").Value = -166667622
$ws.Range("N67").Value = -10591.5
$ws.Range("H87").Value = 63085.25
$ws.Range("J87").Value = 63085.25
$ws.Range("L87").Value = 63085.25
$ws.Range("N87").Value = -65581.25
$ws.Range("H90").Value = 63085.25
$ws.Range("J90").Value = 63085.25
$ws.Range("L90").Value = 189255.75
$ws.Range("N90").Value = -201735.75
$ws.Range("H106").Value = 11294
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null
$ws.Range("H131").Value = 2285.7778
$ws.Range("I131").Value = 2285.7778
$ws.Range("K131").Value = 6857.3334
$ws.Range("M131").Value = -1817.3334
$ws.Range("H132").Value = 223954.27
$ws.Range("I132").Value = 255845.08
$ws.Range("J132").Value = 8691.25
$ws.Range("K132").Value = 767535.24
$ws.Range("L132").Value = 26073.75
$ws.Range("M132").Value = -765005.24
$ws.Range("N132").Value = -31133.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 13672.5
$ws.Range("J44").Value = 13672.5
$ws.Range("L44").Value = 13672.5
$ws.Range("N44").Value = -14648.5
$ws.Range("H80").Value = 27222
$ws.Range("J80").Value = 44444
$ws.Range("L80").Value = 44444
$ws.Range("N80").Value = -46440
$ws.Range("H83").Value = 27222
$ws.Range("J83").Value = 44444
$ws.Range("L83").Value = 133332
$ws.Range("N83").Value = -143316
$ws.Range("H122").Value = 4366.64
$ws.Range("I122").Value = 3703.5264
$ws.Range("K122").Value = 11110.5792
$ws.Range("M122").Value = -8660.5792
$ws.Range("H132").Value = 578143.1
$ws.Range("I132").Value = 676475.1
$ws.Range("J132").Value = 86483.27
$ws.Range("K132").Value = 2029425.3
$ws.Range("L132").Value = 259449.81
$ws.Range("M132").Value = -2026895.3
$ws.Range("N132").Value = -264509.81
$ws.Range("H139").Value = 65902
$ws.Range("J139").Value = 65902
$ws.Range("L139").Value = 65902
$ws.Range("N139").Value = -76182
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6355.143
$ws.Range("I86").Value = 3620.3333
$ws.Range("K86").Value = 3620.3333
$ws.Range("M86").Value = -2497.3333
$ws.Range("H89").Value = 6355.143
$ws.Range("I89").Value = 3620.3333
$ws.Range("K89").Value = 18101.6665
$ws.Range("M89").Value = -12485.6665
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 37643.668
$ws.Range("J9").Value = 37643.668
$ws.Range("L9").Value = 37643.668
$ws.Range("N9").Value = -37979.668
$ws.Range("H105").Value = 500004000
$ws.Range("I105").Value = 1000000000
$ws.Range("K105").Value = 1000000000
$ws.Range("M105").Value = -999998253
$ws.Range("H122").Value = 5994.4707
$ws.Range("I122").Value = 1914.7693
$ws.Range("K122").Value = 5744.3079
$ws.Range("M122").Value = -3294.3079
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 16671685
$ws.Range("I131").Value = 20004618
$ws.Range("J131").Value = 7015.5
$ws.Range("K131").Value = 60013854
$ws.Range("L131").Value = 21046.5
$ws.Range("M131").Value = -60008814
$ws.Range("N131").Value = -31126.5
$ws.Range("H136").Value = 41671948
$ws.Range("I136").Value = 18523784
$ws.Range("J136").Value = 111116450
$ws.Range("K136").Value = 55571352
$ws.Range("L136").Value = 333349350
$ws.Range("M136").Value = -55566252
$ws.Range("N136").Value = -333359550
$ws.Range("H140").Value = 28848384
$ws.Range("I140").Value = 35715812
$ws.Range("J140").Value = 5194
$ws.Range("K140").Value = 107147436
$ws.Range("L140").Value = 15582
$ws.Range("M140").Value = -107142256
$ws.Range("N140").Value = -25942
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 30010.25
$ws.Range("I58").Value = 30010.25
$ws.Range("K58").Value = 30010.25
$ws.Range("M58").Value = -29733.25
$ws.Range("H70").Value = 8129
$ws.Range("I70").Value = 8129
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 8129
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -7859
$ws.Range("N70").Value = $null
$ws.Range("H73").Value = 8129
$ws.Range("I73").Value = 8129
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 8129
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -7193
$ws.Range("N73").Value = $null
$ws.Range("H102").Value = 1005345.7
$ws.Range("I102").Value = 1670602
$ws.Range("J102").Value = 7461.1665
$ws.Range("K102").Value = 1670602
$ws.Range("L102").Value = 7461.1665
$ws.Range("M102").Value = -1668980
$ws.Range("N102").Value = -10705.1665
$ws.Range("H122").Value = 4728.2593
$ws.Range("I122").Value = 3431.4736
$ws.Range("K122").Value = 10294.4208
$ws.Range("M122").Value = -7844.4208
$ws.Range("H132").Value = 7012.6665
$ws.Range("I132").Value = 7286.75
$ws.Range("J132").Value = 5916.3335
$ws.Range("K132").Value = 21860.25
$ws.Range("L132").Value = 17749.0005
$ws.Range("M132").Value = -19330.25
$ws.Range("N132").Value = -22809.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 880
$ws.Range("I22").Value = 796.6667
$ws.Range("J22").Value = 963.3333
$ws.Range("K22").Value = 796.6667
$ws.Range("L22").Value = 963.3333
$ws.Range("M22").Value = -501.6667
$ws.Range("N22").Value = -1553.3333
$ws.Range("H27").Value = 880
$ws.Range("I27").Value = 796.6667
$ws.Range("J27").Value = 963.3333
$ws.Range("K27").Value = 796.6667
$ws.Range("L27").Value = 963.3333
$ws.Range("M27").Value = -689.6667
$ws.Range("N27").Value = -1177.3333
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8237.6
$ws.Range("I122").Value = 1189.8667
$ws.Range("K122").Value = 3569.6001
$ws.Range("M122").Value = -1119.6001
$ws.Range("H126").Value = 4042.7
$ws.Range("I126").Value = 2297.4443
$ws.Range("K126").Value = 6892.3329
$ws.Range("M126").Value = -4422.3329
